$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Title: "Previously done work" -> "Previous Work" ---
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Previous Work"

# --- Content placeholder: give it an explicit position/size (matches the ---
# --- inherited layout/master geometry) and fill in the "Previous Work" ---
# --- bullet content. ---
$body = $s.Shapes.Item(2)
$body.Left = 66
$body.Top = 143.75
$body.Width = 828
$body.Height = 342.6251

$tr = $body.TextFrame.TextRange

$tr.Text = "Worked on implementing an AIM System for a "
$null = $tr.InsertAfter("previous class")
$null = $tr.InsertAfter("`r" + [char]0x201C + "A Multiagent Approach to Autonomous Intersection Management" + [char]0x201D + " by Kurt ")
$null = $tr.InsertAfter("Dresner")
$null = $tr.InsertAfter(" and Peter Stone")
$null = $tr.InsertAfter("`rReservation system using first-come-first-serve policy")
$null = $tr.InsertAfter("`rUse of a grid of " + [char]0x201C + "reservation" + [char]0x201D + " tiles")
$null = $tr.InsertAfter("`rDealing with acceleration in the intersection. Consider trajectories where the vehicle accelerates to max velocity, and maintains current velocity")
$null = $tr.InsertAfter("`rForce a minimum velocity")
$null = $tr.InsertAfter("`rAbility to change the policy")

# Second-level bullets (OOXML lvl="1") for the last five paragraphs.
for ($i = 3; $i -le 7; $i++) {
    $tr.Paragraphs($i).IndentLevel = 2
}
